$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "COB"
$ws.Range("C21").Value = "entity [BFO:0000001]"
$ws.Range("D21").Value = "planned process [COB:0000082]"
$ws.Range("E21").Value = "all"
